$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 8-10 previously had no "stok" (column F) value; fill them in.
$ws.Cells.Item(8,6).Value = "var"
$ws.Cells.Item(9,6).Value = "Var"
$ws.Cells.Item(10,6).Value = "Var"

# Append the new "Gömlek" category products as rows 11-30.
$r = 11
$ws.Cells.Item($r,1).Value = "Kot Gömlek Bej"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "BEJ.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 12
$ws.Cells.Item($r,1).Value = "Kot Gömlek Bordo"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "BORDOKOTGMLK.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 13
$ws.Cells.Item($r,1).Value = "Kot Gömlek Denim Blue"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "DENİMBLUEKOTGMLK.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 14
$ws.Cells.Item($r,1).Value = "Kot Gömlek Fıstık Yeşili"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "FISTIKYEŞİLİ.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 15
$ws.Cells.Item($r,1).Value = "Kot Gömlek Füme"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "FÜMEKOTGMLK.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 16
$ws.Cells.Item($r,1).Value = "Kot Gömlek Kahverengi"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "KAHVERENGİ.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 17
$ws.Cells.Item($r,1).Value = "Kot Gömlek Kraliyet Mavisi"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "KRALİYET MAVİSİ.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 18
$ws.Cells.Item($r,1).Value = "Kot Gömlek Peygamber Çiçeği"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "PEYGAMBERÇİÇEĞİ.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 19
$ws.Cells.Item($r,1).Value = "Kot Gömlek Siyah"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "SİYAHKOTGMLK.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 20
$ws.Cells.Item($r,1).Value = "Kot Gömlek Mavi"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "MAVİKOTGMLK.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 21
$ws.Cells.Item($r,1).Value = "Kot Gömlek Taş"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "TAŞ.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 22
$ws.Cells.Item($r,1).Value = "Kot Gömlek Yeşilin 51. Tonu"
$ws.Cells.Item($r,2).Value = "300 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "YEŞİLİN51.TONU.jpg"
$ws.Cells.Item($r,5).Value = "%85 pamuk, %12 polyester ve %3 spandex karışımından oluşan materyali ile konforlu bir deneyim sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 23
$ws.Cells.Item($r,1).Value = "Kot Gömlek Ceket Denim Blue"
$ws.Cells.Item($r,2).Value = "400 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "GÖMLEK CEKETDENİMBLUE.jpg"
$ws.Cells.Item($r,5).Value = "%100 pamuk materyali ile üretilmiş, Çıt çıt kapama şekliyle güvenli bir kullanım sağlar.İki tarz tek parçada;Hem gömlek hem ceket.XS-S /  M-L /  XL-2XL Beden seçeneği Mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 24
$ws.Cells.Item($r,1).Value = "Kot Gömlek Ceket Mavi"
$ws.Cells.Item($r,2).Value = "400 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "GÖMLEK CEKETMAVİ.jpg"
$ws.Cells.Item($r,5).Value = "%100 pamuk materyali ile üretilmiş, Çıt çıt kapama şekliyle güvenli bir kullanım sağlar.İki tarz tek parçada;Hem gömlek hem ceket.XS-S /  M-L /  XL-2XL Beden seçeneği Mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 25
$ws.Cells.Item($r,1).Value = "Kot Gömlek Ceket Kar Yıkama"
$ws.Cells.Item($r,2).Value = "400 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "GÖMLEKCEKETBUZMAVİSİ.jpg"
$ws.Cells.Item($r,5).Value = "%100 pamuk materyali ile üretilmiş, Çıt çıt kapama şekliyle güvenli bir kullanım sağlar.İki tarz tek parçada;Hem gömlek hem ceket.XS-S /  M-L /  XL-2XL Beden seçeneği Mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 26
$ws.Cells.Item($r,1).Value = "Eşref Gömlek Açık Mavi"
$ws.Cells.Item($r,2).Value = "260 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "EŞREFAÇIKMAVİ.jpg"
$ws.Cells.Item($r,5).Value = "Polyester materyali sayesinde hafifliğiyle konforlu bir kullanım sunarken dayanıklılığından da taviz vermez.Uzun kollu oluşu soğuk havalarda ekstra koruma sağlarken çizgili deseniyle trendleri takip eder.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 27
$ws.Cells.Item($r,1).Value = "Eşref Gömlek Bej"
$ws.Cells.Item($r,2).Value = "260 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "EŞREFBEJ.jpg"
$ws.Cells.Item($r,5).Value = "Polyester materyali sayesinde hafifliğiyle konforlu bir kullanım sunarken dayanıklılığından da taviz vermez.Uzun kollu oluşu soğuk havalarda ekstra koruma sağlarken çizgili deseniyle trendleri takip eder.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 28
$ws.Cells.Item($r,1).Value = "Eşref Gömlek Beyaz"
$ws.Cells.Item($r,2).Value = "260 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "EŞREFBEYAZ.jpg"
$ws.Cells.Item($r,5).Value = "Polyester materyali sayesinde hafifliğiyle konforlu bir kullanım sunarken dayanıklılığından da taviz vermez.Uzun kollu oluşu soğuk havalarda ekstra koruma sağlarken çizgili deseniyle trendleri takip eder.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 29
$ws.Cells.Item($r,1).Value = "Eşref Gömlek Koyu Mavi"
$ws.Cells.Item($r,2).Value = "260 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "EŞREFKOYUMAVİ.jpg"
$ws.Cells.Item($r,5).Value = "Polyester materyali sayesinde hafifliğiyle konforlu bir kullanım sunarken dayanıklılığından da taviz vermez.Uzun kollu oluşu soğuk havalarda ekstra koruma sağlarken çizgili deseniyle trendleri takip eder.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

$r = 30
$ws.Cells.Item($r,1).Value = "Eşref Gömlek Siyah"
$ws.Cells.Item($r,2).Value = "260 Tl"
$ws.Cells.Item($r,3).Value = "Gömlek"
$ws.Cells.Item($r,4).Value = "EŞREFSİYAH.jpg"
$ws.Cells.Item($r,5).Value = "Polyester materyali sayesinde hafifliğiyle konforlu bir kullanım sunarken dayanıklılığından da taviz vermez.Uzun kollu oluşu soğuk havalarda ekstra koruma sağlarken çizgili deseniyle trendleri takip eder.S-M-L-XL-2XL Beden seçeneği mevcuttur."
$ws.Cells.Item($r,6).Value = "Var"

# E26:E30 use a dedicated style (new cellXfs entry with Arial/666666 font).
$r2 = $ws.Range("E26:E30")
$r2.Font.Name = "Arial"
$r2.Font.Size = 11
$r2.Font.Color = 6710886

# Match the final selection recorded in the workbook.
$ws.Range("I32").Select()

